$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 5: update title and link
$ws.Range("D5").Value = "선형 시불변(LTI) 시스템"
$ws.Range("E5").Value = "https://angeloyeo.github.io/2022/01/11/LTI_system.html"

# Row 26: update title only
$ws.Range("D26").Value = "ai plus(est soft)"

# Row 28: update title and link
$ws.Range("D28").Value = "[Manipulator] Dynamics"
$ws.Range("E28").Value = "https://ropiens.tistory.com/175"
